$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "TextBox 83" shape (id 84) that holds the Conclusion / Future
# Work paragraphs on the poster.
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Id -eq 84) {
        $shp = $candidate
        break
    }
}

$tr = $shp.TextFrame.TextRange
$full = $tr.Text

# The shape auto-fits its height to the text (<a:spAutoFit/>); the source
# edit only touched the wording, not the box geometry, so remember the
# original height and restore it once the text has been corrected.
$origHeight = $shp.Height

# --- Fix 1: merge the "voting scheme." run back into the sentence that -----
# --- precedes it, in the "Future Work" paragraph.                      -----
$idx = $full.IndexOf("In the future, there are many further improvements")
$paraLen = $full.Length - $idx
$fixFuture = $tr.Characters($idx + 1, $paraLen)
$fixFuture.Text = "In the future, there are many further improvements we can try.  For example, it is possible to condition the end prediction on the start prediction.  It might also be useful to take multiple models and ensemble them using a voting scheme."

# --- Fix 2: correct the "character leve CNN" typo to "character level CNN", -
# --- moving the word "character" into the same run as "level".        -----
$idx2 = $full.IndexOf("character leve ")
$len2 = "character leve ".Length
$fixTypo = $tr.Characters($idx2 + 1, $len2)
$fixTypo.Text = "character level "

# Restore the original box height now that the text has been patched up.
$shp.Height = $origHeight
